$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# Update column F ("想去人数") values on both the "展览" and "全部类型" sheets.
# The same underlying events appear on both sheets (just at different row
# offsets on "全部类型" because it interleaves rows from other sheets), so
# each update is applied twice, once per sheet/row pair.
$ws1.Range("F6").Value = 571
$ws4.Range("F6").Value = 571
$ws1.Range("F7").Value = 334
$ws4.Range("F7").Value = 334
$ws1.Range("F9").Value = 210
$ws4.Range("F9").Value = 210
$ws1.Range("F10").Value = 220
$ws4.Range("F10").Value = 220
$ws1.Range("F13").Value = 253
$ws4.Range("F13").Value = 253
$ws1.Range("F14").Value = 1042
$ws4.Range("F14").Value = 1042
$ws1.Range("F15").Value = 1380
$ws4.Range("F15").Value = 1380
$ws1.Range("F16").Value = 575
$ws4.Range("F16").Value = 575
$ws1.Range("F17").Value = 94
$ws4.Range("F17").Value = 94
$ws1.Range("F18").Value = 726
$ws4.Range("F18").Value = 726
$ws1.Range("F20").Value = 112
$ws4.Range("F20").Value = 112
$ws1.Range("F21").Value = 102
$ws4.Range("F21").Value = 102
$ws1.Range("F22").Value = 381
$ws4.Range("F22").Value = 381
$ws1.Range("F23").Value = 1248
$ws4.Range("F23").Value = 1248
$ws1.Range("F24").Value = 87
$ws4.Range("F24").Value = 87
$ws1.Range("F25").Value = 58
$ws4.Range("F25").Value = 58
$ws1.Range("F26").Value = 252
$ws4.Range("F26").Value = 252
$ws1.Range("F27").Value = 5212
$ws4.Range("F28").Value = 5212
$ws1.Range("F28").Value = 625
$ws4.Range("F29").Value = 625
$ws1.Range("F30").Value = 178
$ws4.Range("F32").Value = 178
$ws1.Range("F31").Value = 5327
$ws4.Range("F34").Value = 5327
$ws1.Range("F33").Value = 107
$ws4.Range("F36").Value = 107
$ws1.Range("F36").Value = 13645
$ws4.Range("F39").Value = 13645
$ws1.Range("F37").Value = 1392
$ws4.Range("F40").Value = 1392
$ws1.Range("F38").Value = 173
$ws4.Range("F41").Value = 173
$ws1.Range("F39").Value = 69
$ws4.Range("F42").Value = 69
$ws1.Range("F40").Value = 80
$ws4.Range("F43").Value = 80
$ws1.Range("F41").Value = 360
$ws4.Range("F44").Value = 360
$ws1.Range("F42").Value = 508
$ws4.Range("F45").Value = 508
$ws1.Range("F43").Value = 4122
$ws4.Range("F46").Value = 4122
$ws1.Range("F44").Value = 60
$ws4.Range("F47").Value = 60
$ws1.Range("F45").Value = 341
$ws4.Range("F48").Value = 341
$ws1.Range("F46").Value = 107
$ws4.Range("F49").Value = 107

# Row 7 event is sold out: G column switches from a numeric price to the
# text "已售罄" (sold out).
$ws1.Range("G7").Value = "已售罄"
$ws4.Range("G7").Value = "已售罄"

